$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(1)
try {
    $sh = $layout.Shapes.AddPlaceholder(2)
    Write-Output ("Added shape id=" + $sh.Id + " name=" + $sh.Name)
} catch {
    Write-Output ("Error: " + $_.Exception.Message)
}
Write-Output ("Layout Shapes Count now: " + $layout.Shapes.Count)
